# Add three new input-estimate rows ("soil_stones", "landuse_increase_soc",
# "soc") right above the existing "co2_per_egg" row (row 73) on the single
# worksheet, shifting the rows below down by three. This backs the new
# "SOC increase due to landscape change" calculation that feeds the CO2
# certificate estimate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at row 73 - everything from the old row 73
# downward (co2_per_egg ... summerbarley_price) shifts to rows 76-92.
$ws.Rows("73:75").Insert()

# Row 73: soil_stones
$ws.Range("A73").Value = "soil_stones"
$ws.Range("B73").Value = 0.35
$ws.Range("C73").Value = "NA"
$ws.Range("D73").Value = 0.45
$ws.Range("E73").Value = "posnorm"
$ws.Range("F73").Value = "Part of stones in soil in %"

# Row 74: landuse_increase_soc
$ws.Range("A74").Value = "landuse_increase_soc"
$ws.Range("B74").Value = 0.4
$ws.Range("C74").Value = "NA"
$ws.Range("D74").Value = 0.6
$ws.Range("E74").Value = "posnorm"

# Row 75: soc
$ws.Range("A75").Value = "soc"
$ws.Range("B75").Value = 33000
$ws.Range("C75").Value = "NA"
$ws.Range("D75").Value = 41000
$ws.Range("E75").Value = "posnorm"
$ws.Range("F75").Value = "SOC in area of interest in kg"

# Label for row 74 is filled in last to match the original authoring order
# (and therefore the shared-string table order) of the source workbook.
$ws.Range("F74").Value = "SOC increase"

# Scroll the view down to the newly added rows and leave the selection
# where the author left it after the edit.
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G77").Select()
